# Atualizado por script em 13-11-2023 22:16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the match-detail columns (F:V) between row 146 and row 150 ---
# Columns A-E (Indice, pais, torneio, temporada, data_partida) stay put; only
# the actual game data (home/away teams, scores, odds, timestamps, url) swap.
$row146 = @{}
$row150 = @{}
foreach ($col in 6..22) {
    $row146[$col] = $ws.Cells.Item(146, $col).Value()
    $row150[$col] = $ws.Cells.Item(150, $col).Value()
}
foreach ($col in 6..22) {
    $ws.Cells.Item(146, $col).Value = $row150[$col]
    $ws.Cells.Item(150, $col).Value = $row146[$col]
}

# --- Step 2: append new row 151 with the Montana vs Dobrudzha match ---
# Clone formatting from row 150 (A column: bold/centered/bordered "Indice"
# style; E column: datetime number format) before filling in values.
$ws.Cells.Item(150, 1).Copy() | Out-Null
$ws.Cells.Item(151, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(150, 5).Copy() | Out-Null
$ws.Cells.Item(151, 5).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(151, 1).Value = 150
$ws.Cells.Item(151, 2).Value = "bulgaria"
$ws.Cells.Item(151, 3).Value = "vtora-liga"
$ws.Cells.Item(151, 4).Value = "2023-2024"
$ws.Cells.Item(151, 5).Value = 45243.6875
$ws.Cells.Item(151, 6).Value = "Montana"
$ws.Cells.Item(151, 7).Value = 2
$ws.Cells.Item(151, 8).Value = "Dobrudzha"
$ws.Cells.Item(151, 9).Value = 2
$ws.Cells.Item(151, 10).Value = 2.88
$ws.Cells.Item(151, 11).Value = "13/11/2023 05:42"
$ws.Cells.Item(151, 12).Value = 2.14
$ws.Cells.Item(151, 13).Value = "13/11/2023 16:29"
$ws.Cells.Item(151, 14).Value = 2.8
$ws.Cells.Item(151, 15).Value = "13/11/2023 05:42"
$ws.Cells.Item(151, 16).Value = 3.01
$ws.Cells.Item(151, 17).Value = "13/11/2023 16:29"
$ws.Cells.Item(151, 18).Value = 2.46
$ws.Cells.Item(151, 19).Value = "13/11/2023 05:42"
$ws.Cells.Item(151, 20).Value = 3.26
$ws.Cells.Item(151, 21).Value = "13/11/2023 16:29"
$ws.Cells.Item(151, 22).Value = "https://www.betexplorer.com/football/bulgaria/vtora-liga/montana-dobrudzha/hxMgfnnJ/"
